# Add columns with count of topics and modules in the course,
# and remove the now-unused per-profession placeholder sheets.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove all sheets except "Summary" (the first sheet) ---
for ($i = $wb.Worksheets.Count; $i -ge 2; $i--) {
    $wb.Worksheets.Item($i).Delete() | Out-Null
}

$ws = $wb.Worksheets.Item(1)

# --- Set widths for the two new columns (D, E) ---
# Excel's ColumnWidth setter adds a fixed padding offset (~0.8333) when
# translating "characters" into the stored OOXML width, so subtract it
# here to land on the exact target widths of 9 and 8.
$ws.Columns.Item(4).ColumnWidth = 9 - 0.8333333333333333
$ws.Columns.Item(5).ColumnWidth = 8 - 0.8333333333333333

# --- Copy formatting from existing header/data cells onto the new columns ---
# Header row (style matches C1, s="2")
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null

# Data rows (style matches A2, s="3")
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D2:E11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Header text ---
$ws.Range("D1").Value = "Modules"
$ws.Range("E1").Value = "Topics"

# --- Data: number of modules and topics per course ---
$modules = @(11, 19, 19, 12, 9, 20, 25, 20, 16, 8)
$topics = @(90, 148, 185, 40, 64, 138, 182, 192, 166, 45)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $modules[$i]
    $ws.Cells.Item($row, 5).Value = $topics[$i]
}
